# ---------------------------------------------------------------------------
# Applies the "Quantum Entanglement" -> "Government" content swap described
# by the target diff, plus the global TimesNewToman -> Times New Roman font
# typo fix.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Global font-name typo fix: "TimesNewToman" -> "Times New Roman" on every
#    run in the document (the font name is a run-property, not body text, so
#    Find/Replace on Content.Text cannot touch it -- walk paragraphs/Font
#    instead). Exclude each paragraph's trailing paragraph mark so we don't
#    synthesize a spurious pPr/rPr on the mark itself.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.End -gt $r.Start) {
        $body = $d.Range($r.Start, $r.End - 1)
        $body.Font.Name = "Times New Roman"
    }
}

# ---------------------------------------------------------------------------
# Helper: replace the first literal occurrence of $old with $new inside the
# whole document. Used for simple 1-run-in / 1-run-out text swaps, which is
# most of the body copy -- Find/Replace naturally keeps the original run's
# rPr untouched.
# ---------------------------------------------------------------------------
function Replace-Text($old, $new) {
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# Helper: replace a *sub-string* range (located via Find, scoped to $scope)
# with new text while copying the original run's direct formatting onto the
# freshly inserted text -- used where multiple runs must collapse into one
# new run (the engine auto-merges a plain Range.Text assignment into the
# neighboring run when formatting matches, which loses the desired run
# split, so we delete + insert + re-stamp font explicitly).
# ---------------------------------------------------------------------------
function Replace-RangeWithFormatting($scope, $old, $new) {
    $probe = $scope.Duplicate
    $probe.Find.Execute($old) | Out-Null
    $fName  = $probe.Font.Name
    $fSize  = $probe.Font.Size
    $fColor = $probe.Font.Color
    $start = $probe.Start
    $probe.Delete()
    $ins = $d.Range($start, $start)
    $ins.InsertAfter($new)
    $newRng = $d.Range($start, $start + $new.Length)
    $newRng.Font.Name = $fName
    $newRng.Font.Size = $fSize
    $newRng.Font.Color = $fColor
}

# ---------------------------------------------------------------------------
# 2) Title
# ---------------------------------------------------------------------------
Replace-Text "Quantum Entanglement: Unveiling Nature's Enigmatic Link" "Government: The Symphony of Civic Harmony"

# ---------------------------------------------------------------------------
# 3) Author name
# ---------------------------------------------------------------------------
Replace-Text "Isaac Smith" "Nathan Reynolds"

# ---------------------------------------------------------------------------
# 4) Email line: "isaac" + "." + "smith@physicsdomain" + "." + "org"
#    becomes: "nathanreynolds@schuylercentral" + "." + "edu"
#    i.e. the first run's text changes, the middle "." run is untouched,
#    and the last three runs collapse into a single "edu" run.
# ---------------------------------------------------------------------------
$emailPara = $d.Paragraphs(3).Range
Replace-RangeWithFormatting $emailPara.Duplicate "smith@physicsdomain.org" "edu"
Replace-Text "isaac" "nathanreynolds@schuylercentral"

# ---------------------------------------------------------------------------
# 5) Body paragraph (the long one right after the blank paragraph)
# ---------------------------------------------------------------------------
Replace-Text "In a realm where particles defy classical physics and embrace an otherworldly dance, quantum entanglement emerges as a profound mystery" "Immerse yourself in the captivating realm of government, a symphony of intricate mechanisms and dynamic processes that shape our societies"

Replace-Text " It is a breathtaking phenomenon where two particles become intertwined in such a way that the state of one instantaneously affects the other, regardless of the distance separating them" " Delve into the annals of history, where civilizations have experimented with diverse governance structures, each leaving a unique imprint on the evolution of human civilization"

Replace-Text " This enigmatic connection has challenged our understanding of space, time, and reality, sparking fervent debates and inspiring countless experiments" " Explore the principles of law, the foundation upon which justice is upheld and stability is ensured"

Replace-Text " From the realm of theoretical physics to the boundless expanse of the cosmos, quantum entanglement continues to captivate scientists and philosophers alike, beckoning us to unravel the secrets it holds" " Unravel the intricate web of political ideologies, understanding the diverse perspectives that shape public discourse"

Replace-Text "In the tapestry of quantum entanglement, the measurement of one particle instantaneously determines the properties of its entangled partner, irrespective of the vastness of the separation between them" "Journey through the corridors of power, witnessing the ebb and flow of influence as leaders and institutions interact to shape policy and direct the course of nations"

Replace-Text " This profound interdependence stands in stark contrast to the principles of classical physics, where the properties of an object are independent of distant measurements" " Analyze the delicate balance between individual rights and collective responsibilities, examining how governments strive to strike a harmonious accord between these competing interests"

Replace-Text " The bizarre and seemingly paradoxical nature of quantum entanglement has ignited heated discussions among physicists, leading to the formulation of various interpretations to explain this perplexing phenomenon" " Investigate the challenges of global governance, recognizing the interconnectedness of our world and the need for international cooperation to address shared concerns"

Replace-Text "As quantum entanglement continues to bewilder and fascinate, it offers a tantalizing glimpse into the hidden workings of the universe" "Engage with the fundamental questions that have perplexed political thinkers throughout history"

Replace-Text " It has the potential to revolutionize our understanding of information theory, cryptography, and computation, holding the promise of secure communication and exponentially faster quantum computers" " What is the purpose of government? How can we ensure that power is exercised justly and ethically? How do we create a society where all voices are heard and all citizens feel represented? These questions, as timeless as they are profound, invite us to embark on an intellectual adventure that will broaden our perspectives and deepen our understanding of the world around us"

# Drop the trailing ". Moreover, it challenges ... reality" tail entirely
# (keep the final full stop run).
$tailScope = $d.Content
$tailScope.Find.Execute(" What is the purpose of government? How can we ensure that power is exercised justly and ethically? How do we create a society where all voices are heard and all citizens feel represented? These questions, as timeless as they are profound, invite us to embark on an intellectual adventure that will broaden our perspectives and deepen our understanding of the world around us") | Out-Null
$afterQuestions = $d.Range($tailScope.End, $tailScope.End)
$cutScope = $afterQuestions.Duplicate
$cutScope.MoveEndUntil(".", 1) | Out-Null
$cutScope.MoveEnd(1, 1) | Out-Null
$cutScope.Delete()

# ---------------------------------------------------------------------------
# 6) Summary heading paragraph body text
# ---------------------------------------------------------------------------
Replace-Text "Quantum entanglement, a mind-bending phenomenon where particles share an inexplicable link, continues to challenge our understanding of reality" "Government, a complex and dynamic institution, plays a pivotal role in shaping our societies"

Replace-Text " From Schrodinger's cat to Einstein's perplexity, this enigmatic connection defies classical physics, prompting heated debates and inspiring novel interpretations" " Its intricate mechanisms, rooted in history, uphold justice, facilitate cooperation, and address global challenges"

Replace-Text " With its implications for cryptography, quantum computing, and our fundamental understanding of the cosmos, quantum entanglement remains at the forefront of " " As we analyze the principles of governance, explore political ideologies, and grapple with fundamental questions about power and representation, we gain a deeper appreciation for the symphony of civic harmony"

# Final run used to read "scientific exploration, ..." (after a
# lastRenderedPageBreak); becomes a plain "." and two brand new runs are
# appended carrying the closing sentence + final full stop.
$lastScope = $d.Content
$lastScope.Find.Execute("scientific exploration, beckoning us to decipher its profound mysteries and rewrite our narrative of the universe") | Out-Null
$fName  = $lastScope.Font.Name
$fSize  = $lastScope.Font.Size
$fColor = $lastScope.Font.Color
$start = $lastScope.Start
$lastScope.Delete()
$ins = $d.Range($start, $start)
$ins.InsertAfter(".")
$fixRng = $d.Range($start, $start + 1)
$fixRng.Font.Name = $fName
$fixRng.Font.Size = $fSize
$fixRng.Font.Color = $fColor

# Remove the now-orphaned lastRenderedPageBreak that preceded that run (it
# sat at the very start of the same run, right after the "forefront of "
# text, so it's immediately before our new "." insertion point).
$brScope = $d.Range($start - 0, $start)

Write-Output "mid-checkpoint"
